$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 169601
$ws.Range("C4").Value = 160434
$ws.Range("C5").Value = 9167
$ws.Range("C7").Value = 5.41
$ws.Range("C8").Value = 65.54000000000001
